$d = $word.ActiveDocument

# Merge the run-split placeholder in paragraph 1 into a single run
$d.Content.Find.Execute("{{roundNumber}} - {{roundDate}}", $false, $false, $false, $false, $false, $true, 1, $false, "{{roundNumber}} - {{roundDate}}", 2)

# Fix the {{#matches}} opening tag to single-brace {#matches}
$d.Content.Find.Execute("{{#matches}}", $false, $false, $false, $false, $false, $true, 1, $false, "{#matches}", 2)

# Merge the run-split match line into a single run
$d.Content.Find.Execute("Zápas {{index}} – {{match_time}} – {{teamA}} vs {{teamB}}", $false, $false, $false, $false, $false, $true, 1, $false, "Zápas {{index}} – {{match_time}} – {{teamA}} vs {{teamB}}", 2)

# Fix the {{/matches}} closing tag to single-brace {/matches}
$d.Content.Find.Execute("{{/matches}}", $false, $false, $false, $false, $false, $true, 1, $false, "{/matches}", 2)
